$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. "18" -> "17" in the first paragraph (task/problem number).
#    The digit "8" only ever occurs as literal document text in that single
#    run, so a plain Find/Replace on the document body text is unambiguous.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("8", $false, $false, $false, $false, $false, $true, 1, $false, "7", 2)

# ---------------------------------------------------------------------------
# 2. First picture ("Рисунок 8" / ur_36_25, InlineShapes(1)):
#      - widen the crop window on the left/right of the source rectangle
#      - grow the displayed width from 153.75pt to 154.5pt (height untouched)
#    PictureFormat.Crop* is expressed in points against the picture's native
#    (uncropped) pixel size; ScaleWidth likewise scales width only, leaving
#    Height alone (unlike the Width property, which preserves the shape's
#    aspect ratio and would drag Height along with it).
# ---------------------------------------------------------------------------
$pic1 = $d.InlineShapes.Item(1)
$pic1.PictureFormat.CropLeft  = 22.24908
$pic1.PictureFormat.CropRight = 303.41871
$pic1.ScaleWidth = 27.24867724867725

# ---------------------------------------------------------------------------
# 3. Second picture ("Рисунок 9" / slide_30, InlineShapes(2)):
#      - shift the crop window further right (more cropped on the left,
#        less on the right)
#      - grow the displayed width from 161.25pt to 180pt (height untouched)
# ---------------------------------------------------------------------------
$pic2 = $d.InlineShapes.Item(2)
$pic2.PictureFormat.CropLeft  = 458.202
$pic2.PictureFormat.CropRight = 13.038
$pic2.ScaleWidth = 30
